# assignment1.xlsx - "Add files via upload"
# The sales rep for order 3 (row 4) now gets a flat 50% surcharge instead of
# the shipmode-derived lookup formula, and the active selection moves to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the IF(...) lookup formula in E4 with a hard-coded value.
$ws.Range("E4").Value = 0.5

# F4's formula is unchanged; it recalculates automatically from the new E4.

# Update the last-saved selection to E5.
$ws.Range("E5").Select()
